$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2788583.2
$ws.Range("I32").Value = 692.9231
$ws.Range("J32").Value = 5808798
$ws.Range("K32").Value = 692.9231
$ws.Range("L32").Value = 5808798
$ws.Range("M32").Value = -366.9231
$ws.Range("N32").Value = -5809450

$ws.Range("H41").Value = 5437.8096
$ws.Range("I41").Value = 491.1
$ws.Range("J41").Value = 9934.817999999999
$ws.Range("K41").Value = 491.1
$ws.Range("L41").Value = 9934.817999999999
$ws.Range("M41").Value = -51.10000000000002
$ws.Range("N41").Value = -10814.818

$ws.Range("H80").Value = 6677
$ws.Range("I80").Value = 3781
$ws.Range("J80").Value = 9061.941000000001
$ws.Range("K80").Value = 11343
$ws.Range("L80").Value = 27185.823
$ws.Range("M80").Value = -10345
$ws.Range("N80").Value = -29181.823

$ws.Range("H83").Value = 6677
$ws.Range("I83").Value = 3781
$ws.Range("J83").Value = 9061.941000000001
$ws.Range("K83").Value = 34029
$ws.Range("L83").Value = 81557.46900000001
$ws.Range("M83").Value = -29037
$ws.Range("N83").Value = -91541.46900000001

$ws.Range("H86").Value = 5531.952
$ws.Range("J86").Value = 6706.5
$ws.Range("L86").Value = 6706.5
$ws.Range("N86").Value = -8952.5

$ws.Range("H89").Value = 5531.952
$ws.Range("J89").Value = 6706.5
$ws.Range("L89").Value = 33532.5
$ws.Range("N89").Value = -44764.5

$ws.Range("H137").Value = 43160.92
$ws.Range("I137").Value = 60693.723
$ws.Range("J137").Value = 3712.125
$ws.Range("K137").Value = 182081.169
$ws.Range("L137").Value = 11136.375
$ws.Range("M137").Value = -179531.169
$ws.Range("N137").Value = -16236.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2165.0715
$ws.Range("I45").Value = 1876.375
$ws.Range("J45").Value = 2550
$ws.Range("K45").Value = 1876.375
$ws.Range("L45").Value = 2550
$ws.Range("M45").Value = -1499.375
$ws.Range("N45").Value = -3304

$ws.Range("H80").Value = 30110
$ws.Range("J80").Value = 30110
$ws.Range("L80").Value = 30110
$ws.Range("N80").Value = -32106

$ws.Range("H83").Value = 30110
$ws.Range("J83").Value = 30110
$ws.Range("L83").Value = 90330
$ws.Range("N83").Value = -100314

$ws.Range("H122").Value = 3379.8235
$ws.Range("I122").Value = 4129.9
$ws.Range("J122").Value = 2308.2856
$ws.Range("K122").Value = 12389.7
$ws.Range("L122").Value = 6924.8568
$ws.Range("M122").Value = -9939.699999999999
$ws.Range("N122").Value = -11824.8568

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 3262
$ws.Range("I132").Value = 3255.0344
$ws.Range("J132").Value = 3273.8823
$ws.Range("K132").Value = 9765.1032
$ws.Range("L132").Value = 9821.6469
$ws.Range("M132").Value = -7235.1032
$ws.Range("N132").Value = -14881.6469

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5563.654
$ws.Range("I86").Value = 4055.6
$ws.Range("J86").Value = 10590.5
$ws.Range("K86").Value = 4055.6
$ws.Range("L86").Value = 10590.5
$ws.Range("M86").Value = -2932.6
$ws.Range("N86").Value = -12836.5

$ws.Range("H89").Value = 5563.654
$ws.Range("I89").Value = 4055.6
$ws.Range("J89").Value = 10590.5
$ws.Range("K89").Value = 20278
$ws.Range("L89").Value = 52952.5
$ws.Range("M89").Value = -14662
$ws.Range("N89").Value = -64184.5

$ws.Range("H99").Value = 3191.1738
$ws.Range("I99").Value = 3678.842
$ws.Range("J99").Value = 874.75
$ws.Range("K99").Value = 3678.842
$ws.Range("L99").Value = 874.75
$ws.Range("M99").Value = -2180.842
$ws.Range("N99").Value = -3870.75

$ws.Range("H107").Value = 1986.0869
$ws.Range("I107").Value = 2044.6154
$ws.Range("J107").Value = 1910
$ws.Range("K107").Value = 2044.6154
$ws.Range("L107").Value = 1910
$ws.Range("M107").Value = -124.6153999999999
$ws.Range("N107").Value = -5750

$ws.Range("H134").Value = 29498.975
$ws.Range("I134").Value = 49255.24
$ws.Range("K134").Value = 147765.72
$ws.Range("M134").Value = -145230.72

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 624.8333
$ws.Range("I22").Value = 519.6667
$ws.Range("J22").Value = 730
$ws.Range("K22").Value = 519.6667
$ws.Range("L22").Value = 730
$ws.Range("M22").Value = -169.6667
$ws.Range("N22").Value = -1430

$ws.Range("H107").Value = 678.7778
$ws.Range("I107").Value = 323.8
$ws.Range("J107").Value = 1122.5
$ws.Range("K107").Value = 323.8
$ws.Range("L107").Value = 1122.5
$ws.Range("M107").Value = 1596.2
$ws.Range("N107").Value = -4962.5

$ws.Range("H132").Value = 1755.2931
$ws.Range("I132").Value = 1055.7894
$ws.Range("J132").Value = 3084.35
$ws.Range("K132").Value = 3167.3682
$ws.Range("L132").Value = 9253.049999999999
$ws.Range("M132").Value = -637.3681999999999
$ws.Range("N132").Value = -14313.05

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2026
$ws.Range("I62").Value = 998
$ws.Range("J62").Value = 2540
$ws.Range("K62").Value = 2994
$ws.Range("L62").Value = 7620
$ws.Range("M62").Value = -2308
$ws.Range("N62").Value = -8992

$ws.Range("H65").Value = 2026
$ws.Range("I65").Value = 998
$ws.Range("J65").Value = 2540
$ws.Range("K65").Value = 8982
$ws.Range("L65").Value = 22860
$ws.Range("M65").Value = -5550
$ws.Range("N65").Value = -29724

$ws.Range("H108").Value = 2522.2856
$ws.Range("I108").Value = 478
$ws.Range("J108").Value = 3340
$ws.Range("K108").Value = 1434
$ws.Range("L108").Value = 10020
$ws.Range("M108").Value = 1446
$ws.Range("N108").Value = -15780

$ws.Range("H109").Value = 2725
$ws.Range("I109").Value = 933.3333
$ws.Range("J109").Value = 3800
$ws.Range("K109").Value = 2799.9999
$ws.Range("L109").Value = 11400
$ws.Range("M109").Value = -1759.9999
$ws.Range("N109").Value = -13480

$ws.Range("H119").Value = 2583.3333
$ws.Range("I119").Value = 1875
$ws.Range("J119").Value = 4000
$ws.Range("K119").Value = 5625
$ws.Range("L119").Value = 12000
$ws.Range("M119").Value = -787
$ws.Range("N119").Value = -21676

$ws.Range("H131").Value = 914.03174
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 914.03174
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2742.09522
$ws.Range("N131").Value = -12822.09522
$ws.Range("M131").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4621.1816
$ws.Range("I80").Value = 4799.967
$ws.Range("K80").Value = 4799.967
$ws.Range("M80").Value = -3801.967

$ws.Range("H83").Value = 4621.1816
$ws.Range("I83").Value = 4799.967
$ws.Range("K83").Value = 23999.835
$ws.Range("M83").Value = -19007.835

$ws.Range("H132").Value = 3775.0908
$ws.Range("I132").Value = 4142.6665
$ws.Range("J132").Value = 3334
$ws.Range("K132").Value = 12427.9995
$ws.Range("L132").Value = 10002
$ws.Range("M132").Value = -9897.999500000002
$ws.Range("N132").Value = -15062

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 707.3333
$ws.Range("I46").Value = 262
$ws.Range("J46").Value = 930
$ws.Range("K46").Value = 262
$ws.Range("L46").Value = 930
$ws.Range("M46").Value = -74
$ws.Range("N46").Value = -1306

$ws.Range("H68").Value = 3581.182
$ws.Range("I68").Value = 2374
$ws.Range("J68").Value = 3936.2354
$ws.Range("K68").Value = 2374
$ws.Range("L68").Value = 3936.2354
$ws.Range("M68").Value = -1625
$ws.Range("N68").Value = -5434.2354

$ws.Range("H71").Value = 3581.182
$ws.Range("I71").Value = 2374
$ws.Range("J71").Value = 3936.2354
$ws.Range("K71").Value = 11870
$ws.Range("L71").Value = 19681.177
$ws.Range("M71").Value = -8126
$ws.Range("N71").Value = -27169.177

$ws.Range("H93").Value = 1480.4333
$ws.Range("I93").Value = 1386.9546
$ws.Range("J93").Value = 1737.5
$ws.Range("K93").Value = 1386.9546
$ws.Range("L93").Value = 1737.5
$ws.Range("M93").Value = -138.9546
$ws.Range("N93").Value = -4233.5

$ws.Range("H100").Value = 33337112
$ws.Range("I100").Value = 4765.1875
$ws.Range("J100").Value = 71431224
$ws.Range("K100").Value = 4765.1875
$ws.Range("L100").Value = 71431224
$ws.Range("M100").Value = -4224.1875
$ws.Range("N100").Value = -71432306

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 85383.914
$ws.Range("I122").Value = 126825.875
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 380477.625
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -378027.625
$ws.Range("N122").Value = -12400

$ws.Range("I132").Value = 68396.2
$ws.Range("J132").Value = 3216.4211
$ws.Range("K132").Value = 205188.6
$ws.Range("L132").Value = 9649.263300000001
$ws.Range("M132").Value = -202658.6
$ws.Range("N132").Value = -14709.2633

$ws.Range("H136").Value = 41669760
$ws.Range("I136").Value = 62503004
$ws.Range("J136").Value = 20836514
$ws.Range("K136").Value = 187509012
$ws.Range("L136").Value = 62509542
$ws.Range("M136").Value = -187506462
$ws.Range("N136").Value = -62514642
